# tags.xlsx: replace the single "Description" header (F1) with a
# "value" column (F) + a new "description" column (G), and add the
# boolean-ish numeric "value" data for each tag row (F2:F4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("F1").Value = "value"
$ws.Range("G1").Value = "description"

# New "value" data for the existing tag rows
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 1

# Match the author's final UI selection state
$ws.Range("G9").Select() | Out-Null
